# Updated symbol list on Sat Jan 14 17:52:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.99"
$ws.Range("E2").Value = "'6.10%"
$ws.Range("D3").Value = "'32.23"
$ws.Range("E3").Value = "'10.76%"
$ws.Range("D4").Value = "'5.295"
$ws.Range("E4").Value = "'1.92%"
$ws.Range("D5").Value = "'0.07489"
$ws.Range("E5").Value = "'7.52%"
$ws.Range("D6").Value = "'7.856"
$ws.Range("E6").Value = "'5.98%"
$ws.Range("D7").Value = "'3.806"
$ws.Range("E7").Value = "'6.95%"
$ws.Range("D8").Value = "'1.571"
$ws.Range("E8").Value = "'12.83%"
$ws.Range("D9").Value = "'0.9186"
$ws.Range("E9").Value = "'2.39%"
$ws.Range("D10").Value = "'0.01761"
$ws.Range("E10").Value = "'2,620.43%"
$ws.Range("D11").Value = "'0.1699"
$ws.Range("E11").Value = "'6.04%"
$ws.Range("D12").Value = "'0.07825"
$ws.Range("E12").Value = "'2.56%"
$ws.Range("D13").Value = "'0.08021"
$ws.Range("E13").Value = "'4.77%"
$ws.Range("D14").Value = "'0.03001"
$ws.Range("E14").Value = "'2.42%"
$ws.Range("D15").Value = "'0.09911"
$ws.Range("E15").Value = "'10.23%"
$ws.Range("D16").Value = "'0.001493"
$ws.Range("E16").Value = "'-5.72%"
$ws.Range("D17").Value = "'0.04608"
$ws.Range("E17").Value = "'2.01%"
$ws.Range("D18").Value = "'0.006260"
$ws.Range("E18").Value = "'-0.64%"
$ws.Range("D19").Value = "'3.474"
$ws.Range("E19").Value = "'0.39%"
$ws.Range("E20").Value = "'-0.02%"
$ws.Range("D21").Value = "'0.3327"
$ws.Range("E21").Value = "'3.08%"
$ws.Range("D22").Value = "'0.1334"
$ws.Range("E22").Value = "'0.23%"
$ws.Range("D23").Value = "'4.488"
$ws.Range("E23").Value = "'11.92%"
$ws.Range("D24").Value = "'0.1620"
$ws.Range("E24").Value = "'1.42%"
$ws.Range("E25").Value = "'0.80%"
$ws.Range("D26").Value = "'0.004451"
$ws.Range("E26").Value = "'5.15%"
$ws.Range("D27").Value = "'0.0001398"
$ws.Range("E27").Value = "'19.61%"
$ws.Range("D28").Value = "'0.0001746"
$ws.Range("E28").Value = "'6.62%"
$ws.Range("D40").Value = "'0.04537"
$ws.Range("E40").Value = "'5.23%"
$ws.Range("D41").Value = "'0.007178"
$ws.Range("E41").Value = "'3.73%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'8.26%"
$ws.Range("D43").Value = "'0.002217"
$ws.Range("E43").Value = "'6.69%"
$ws.Range("D44").Value = "'0.01283"
$ws.Range("E44").Value = "'11.65%"
$ws.Range("D45").Value = "'0.00006216"
$ws.Range("E45").Value = "'6.51%"
$ws.Range("D46").Value = "'0.7092"
$ws.Range("E46").Value = "'-63.24%"
$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-0.50%"
